$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in B1 (force text so Excel doesn't convert the
# date-shaped string into a date serial number; ClearFormats then drops
# the temporary "@" number-format again so no stray style sticks around)
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2024-12-14"
$ws.Range("B1").ClearFormats()

# Mark cyrus (row 2) as PRESENT with a time
$ws.Range("B2").Value = "PRESENT"
$ws.Range("C2").Value = "14:52:41"

# Update kiefer's (row 3) time
$ws.Range("C3").Value = "14:53:00"

# Mark roche (row 4) as PRESENT with a time
$ws.Range("B4").Value = "PRESENT"
$ws.Range("C4").Value = "14:52:42"
